$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B,C,D,F,G,I,K,L,M,N hold the per-line loading percentages;
# columns E,H,J,O stay at 0 and are left untouched.
$cols = @("B", "C", "D", "F", "G", "I", "K", "L", "M", "N")

# New loading-percent values for rows 2-25 (one row per data point),
# in the same column order as $cols above.
$newValues = @(
    @(15.39256106555644, 5.281295514341418, 7.75811563821173, 38.21005944343406, 3.717395060308168, 31.72469780236165, 13.20996876328252, 10.84794232929675, 16.03704142207095, 23.02341862152726),
    @(15.22459039801095, 5.043646997262314, 7.762640867588784, 38.14602805274244, 3.72013394041565, 31.75301386769192, 13.08526699425582, 10.85604394937295, 16.02128896020184, 23.07429295256215),
    @(15.12460205486797, 4.890689443285682, 7.765457926244703, 38.11492459770506, 3.721904545746547, 31.77570811930995, 13.01152383139821, 10.86261525325722, 16.0145308584872, 23.10742640515188),
    @(15.08469421130342, 4.826627285134463, 7.766615734008594, 38.10432133190169, 3.722648517511951, 31.78628974383145, 12.98221519503958, 10.86569500743624, 16.01251248749491, 23.12140591941487),
    @(15.07811951131141, 4.815886571820998, 7.766808585129652, 38.10268595058893, 3.722773410741981, 31.78812731412706, 12.97739426067564, 10.86623068009512, 16.01222184280842, 23.12375605349744),
    @(15.12406038896298, 4.889832428922744, 7.765473500858071, 38.11477320291876, 3.721914488268115, 31.77584542930981, 13.01112551801749, 10.86265516025104, 16.01450065611857, 23.10761300431298),
    @(15.33402049492961, 5.200837331404946, 7.759668021099543, 38.18628149398253, 3.718321016009737, 31.73335884689218, 13.16640514493287, 10.85040456852665, 16.03100702439325, 23.040566789481),
    @(15.76857376024182, 5.753427968865489, 7.748582756585293, 38.39130470801427, 3.711976335882272, 31.69220446947247, 13.49183822612719, 10.83903400021548, 16.08635399517659, 22.92411375871141),
    @(16.09880482481857, 6.122961669186052, 7.74061127303463, 38.58084624591769, 3.707738100945887, 31.6877167858251, 13.74163510101783, 10.83836607623946, 16.14080885088437, 22.84768302403141),
    @(16.2507776620272, 6.282923867167529, 7.737020311994296, 38.67536158616882, 3.705900885528891, 31.69126865367582, 13.85714431670612, 10.83972341170422, 16.16852124701248, 22.81488810179187),
    @(16.3085238315045, 6.342311993804358, 7.735665433920076, 38.71232799698559, 3.705218155549796, 31.69341737696569, 13.90111576645617, 10.84047541782582, 16.17943263661768, 22.80275292861721),
    @(16.29607924028731, 6.3295746638272, 7.735957013558591, 38.70431462417002, 3.705364617329475, 31.69291887984701, 13.89163609324441, 10.84030288903993, 16.17706420675907, 22.80535384844777),
    @(16.25552482733365, 6.287833612842086, 7.736908747112782, 38.67837937865323, 3.705844457112299, 31.69142932786685, 13.8607574762163, 10.83978051488968, 16.16941060347618, 22.81388405251252),
    @(16.23070818057426, 6.262111143377955, 7.737492350826672, 38.66264584842732, 3.706140061608337, 31.69062157575375, 13.84187230513826, 10.83949151450273, 16.16477672403875, 22.81914596772574),
    @(16.08890359862616, 6.112342559511704, 7.740846649584392, 38.57483489691739, 3.707859987979992, 31.6875971726056, 13.73412070178998, 10.83831072698589, 16.13905648873415, 22.8498659430546),
    @(16.00232039958065, 6.018367368277945, 7.74291334896582, 38.52307801174816, 3.708938307605926, 31.68717408508982, 13.66847114708538, 10.83801123926978, 16.12402730792375, 22.86921698709365),
    @(15.95268788916133, 5.963549807691793, 7.7441053910539, 38.49409041310138, 3.709567077243679, 31.68745726970687, 13.63089012722254, 10.83799539804058, 16.11566006463306, 22.88053299473503),
    @(15.93591366031894, 5.944858569034503, 7.744509572350698, 38.48441043659328, 3.709781438238516, 31.68764359563527, 13.61819774920203, 10.8380169236558, 16.11287482154853, 22.88439632046241),
    @(16.01152035022832, 6.028450537918935, 7.742693001491202, 38.52850684550197, 3.708822634391189, 31.6871646322909, 13.675441424886, 10.83802693767477, 16.12559854532674, 22.86713780796799),
    @(16.26743171450621, 6.3001262696098, 7.736629066880954, 38.68596542868003, 3.7057031647226, 31.69184503985402, 13.86982132702678, 10.83992749664542, 16.17164737306927, 22.81137083140094),
    @(16.43581021136725, 6.470763537886337, 7.732694672200421, 38.79571647758183, 3.703740057203617, 31.69958801015111, 13.99818596204026, 10.8425564665997, 16.204172450257, 22.77657653168949),
    @(16.34585819348575, 6.380328605919582, 7.734791947970428, 38.73652015694861, 3.70478090608993, 31.69502715681372, 13.92956689569928, 10.84102675183056, 16.18659284869031, 22.79499577435769),
    @(16.00736059571671, 6.023894396394883, 7.742792608496181, 38.52605007573202, 3.708874902743154, 31.68716726599073, 13.67228965535755, 10.83801935338539, 16.12488733737214, 22.86807721032488),
    @(15.64887377268307, 5.610247180110944, 7.751550570496299, 38.32895770820753, 3.713618074749729, 31.69881792223926, 13.40176831134879, 10.84075791073374, 16.06894179221491, 22.9540122982508)
)

$startRow = 2
for ($r = 0; $r -lt $newValues.Count; $r++) {
    $rowValues = $newValues[$r]
    for ($c = 0; $c -lt $cols.Count; $c++) {
        $ws.Range($cols[$c] + ($startRow + $r)).Value = $rowValues[$c]
    }
}

